# Update faturamento_diario_lojas.xlsx with new daily figures for days 23-25
# (columns X, Y, Z) and the corresponding updated row totals (column AG).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bibi Cell Mundi
$ws.Range("X2").Value = 15420.59
$ws.Range("Y2").Value = 5091.8
$ws.Range("AG2").Value = 275560.78

# Row 3 - Bibi Cell Vieiralves
$ws.Range("X3").Value = 3054.8
$ws.Range("Y3").Value = 20778.9
$ws.Range("AG3").Value = 167944.3

# Row 4 - Bibi Cell Ponta Negra
$ws.Range("X4").Value = 1998
$ws.Range("Y4").Value = 5791.02
$ws.Range("Z4").Value = 947
$ws.Range("AG4").Value = 82302.32000000001

# Row 5 - Bibi Cell Manauara
$ws.Range("X5").Value = 3346.9
$ws.Range("Y5").Value = 3261
$ws.Range("Z5").Value = 739.9
$ws.Range("AG5").Value = 76326.48

# Row 6 - total
$ws.Range("X6").Value = 23820.29
$ws.Range("Y6").Value = 34922.72
$ws.Range("Z6").Value = 1686.9
$ws.Range("AG6").Value = 602133.88
